# Registration.xlsx edit script
# Implements: new "ChangePassword" sheet, CreateProject text fixes,
# RegisterData row2/3 text-format + D3 text update + pageSetup, selection
# changes on searchProject, and final active-sheet/selection state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "ChangePassword" sheet, then move it to the end
#    (after "CreateProject", which is currently the last sheet).
# ---------------------------------------------------------------------
$wsChange = $wb.Worksheets.Add()
$wsChange.Name = "ChangePassword"
$wsCreate = $wb.Worksheets.Item("CreateProject")
$wsChange.Move($null, $wsCreate)

# Re-fetch handles after the move to avoid stale references.
$wsChange = $wb.Worksheets.Item("ChangePassword")
$wsCreate = $wb.Worksheets.Item("CreateProject")

# ---------------------------------------------------------------------
# 2. CreateProject sheet: fix typo / swap placeholder values
# ---------------------------------------------------------------------
$wsCreate.Range("D1").Value = "Project Type"
$wsCreate.Range("A2").Value = "Vikesh two"
$wsCreate.Range("C2").Value = "aaa"

# ---------------------------------------------------------------------
# 3. Populate the new ChangePassword sheet
# ---------------------------------------------------------------------
$wsChange.Range("A1").Value = "Old Password"
$wsChange.Range("B1").Value = "New Password"
$wsChange.Range("C1").Value = "confirm Password"
$wsChange.Range("A1:C1").Interior.Color = 65535

$wsChange.Range("A2").Value = "Vikesh@1989"
$wsChange.Range("B2").Value = "Test@1989"
$wsChange.Range("C2").Value = "Test@1989"

$wsChange.Columns.Item(1).ColumnWidth = 13.140625
$wsChange.Columns.Item(2).ColumnWidth = 14.140625
$wsChange.Columns.Item(3).ColumnWidth = 16.85546875

# ---------------------------------------------------------------------
# 4. RegisterData sheet: apply text format to rows 2-3, update D3,
#    set page orientation
# ---------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("RegisterData")
$wsReg.Range("A2:H3").NumberFormat = "@"
$wsReg.Range("D3").Value = "All new"
$wsReg.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. searchProject sheet: set a selection (does not stay active tab)
# ---------------------------------------------------------------------
$wsSearch = $wb.Worksheets.Item("searchProject")
$wsSearch.Activate() | Out-Null
$wsSearch.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 6. ChangePassword sheet: set its selection too (not the active tab)
# ---------------------------------------------------------------------
$wsChange.Activate() | Out-Null
$wsChange.Range("H20").Select() | Out-Null

# ---------------------------------------------------------------------
# 7. CreateProject sheet: update selection, no longer the active tab
# ---------------------------------------------------------------------
$wsCreate.Activate() | Out-Null
$wsCreate.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------
# 8. RegisterData becomes the final active sheet / selection
# ---------------------------------------------------------------------
$wsReg.Activate() | Out-Null
$wsReg.Range("G6").Select() | Out-Null
